# Texas COVID-19 Fatality Count by County — "Fully implemented all data for
# the today query": append 3 new daily fatality-count columns (12-30, 12-01,
# 12-02) after the existing last column (JI), and refresh the trailing
# cumulative totals that shifted as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fatalities by County")

# --- add the 3 new date columns --------------------------------------------
# Column JI (269) is the current last column ("Fatalities 11-29"); the new
# columns land in JJ/JK/JL (270/271/272).
$newHeaders = @("Fatalities 12-30", "Fatalities 12-01", "Fatalities 12-02")
$newValues  = @(266, 266, 266)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = 270 + $i
    $headerCell = $ws.Cells.Item(1, $col)
    $headerCell.Value = $newHeaders[$i]
    $ws.Cells.Item(2, $col).Value = $newValues[$i]
}

# Match the look of the rest of the header row (bold/centered/bordered)
# by painting the formatting from the previous header cell (JI1) onto the
# freshly added header cells.
$ws.Range("JI1").Copy() | Out-Null
$ws.Range("JJ1:JL1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- refresh the trailing cumulative totals that moved with the new data ---
# (columns IF .. JI, i.e. 240..269)
$updatedTotals = @{
    240 = 205; 241 = 206; 242 = 209; 243 = 212; 244 = 214; 245 = 215;
    246 = 219; 247 = 221; 248 = 223; 249 = 224; 250 = 225; 251 = 229;
    252 = 231; 253 = 233; 254 = 236; 255 = 239; 256 = 241; 257 = 246;
    258 = 254; 259 = 254; 260 = 254; 261 = 255; 262 = 259; 263 = 261;
    264 = 261; 265 = 263; 266 = 263; 267 = 265; 268 = 265; 269 = 265
}

foreach ($col in $updatedTotals.Keys) {
    $ws.Cells.Item(2, $col).Value = $updatedTotals[$col]
}

# Leave the selection parked on the new last header cell, same as the source
# workbook after the update.
$ws.Range("JH1").Select() | Out-Null
